$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = $true
$ws.Range("D2").Value = 0.2170335732869684
$ws.Range("E2").Value = 0.2170335732869684

# Row 3
$ws.Range("D3").Value = 0.001910844557714726
$ws.Range("E3").Value = 0.001910844557714726

# Row 4
$ws.Range("D4").Value = 0.00278300978396286
$ws.Range("E4").Value = 0.00278300978396286

# Row 5
$ws.Range("D5").Value = 0.7601998640822202
$ws.Range("E5").Value = 0.7601998640822202

# Row 6
$ws.Range("C6").Value = $true
$ws.Range("D6").Value = 0.1136118718498357
$ws.Range("E6").Value = 0.1136118718498357

# Row 7
$ws.Range("D7").Value = 0.02987387769194438
$ws.Range("E7").Value = 0.9701261223080556

# Row 8
$ws.Range("D8").Value = 0.005604820149101408
$ws.Range("E8").Value = 0.9943951798508985

# Row 9
$ws.Range("C9").Value = $false
$ws.Range("D9").Value = 0.03963953114706664
$ws.Range("E9").Value = 0.9603604688529334

# Row 10
$ws.Range("D10").Value = 0.1469175125200937
$ws.Range("E10").Value = 0.8530824874799063
$ws.Range("F10").Value = 1.737625122070312
$ws.Range("G10").Value = 0.4444444444444444

# Row 11
$ws.Range("C11").Value = $true
$ws.Range("D11").Value = 0.2489935434253699
$ws.Range("E11").Value = 0.2489935434253699

# Row 12
$ws.Range("D12").Value = 0.0006165579020269735
$ws.Range("E12").Value = 0.0006165579020269735

# Row 13
$ws.Range("D13").Value = 0.003493858659941289
$ws.Range("E13").Value = 0.003493858659941289

# Row 14
$ws.Range("D14").Value = 0.712600439272515
$ws.Range("E14").Value = 0.712600439272515

# Row 15
$ws.Range("C15").Value = $true
$ws.Range("D15").Value = 0.05203423385868784
$ws.Range("E15").Value = 0.05203423385868784

# Row 16
$ws.Range("D16").Value = 0.01590273578726582
$ws.Range("E16").Value = 0.9840972642127341

# Row 17
$ws.Range("D17").Value = 0.007243261590205338
$ws.Range("E17").Value = 0.9927567384097946

# Row 18
$ws.Range("D18").Value = 0.01908640855893808
$ws.Range("E18").Value = 0.9809135914410619

# Row 19
$ws.Range("D19").Value = 0.05192312716446317
$ws.Range("E19").Value = 0.9480768728355369
$ws.Range("F19").Value = 1.952943801879883
$ws.Range("G19").Value = 0.4444444444444444
